$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 18:35"

# Refresh country statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes)
$ws.Range("B4").Value = 1601035
$ws.Range("C4").Value = 8312
$ws.Range("D4").Value = 371374
$ws.Range("E4").Value = 1134357
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = 95304
$ws.Range("B6").Value = 296033
$ws.Range("C6").Value = 2676
$ws.Range("E6").Value = 160202
$ws.Range("G6").Value = 254
$ws.Range("H6").Value = 19148
$ws.Range("B9").Value = 228006
$ws.Range("C9").Value = 642
$ws.Range("D9").Value = 134560
$ws.Range("E9").Value = 60960
$ws.Range("G9").Value = 156
$ws.Range("H9").Value = 32486
$ws.Range("B11").Value = 178864
$ws.Range("C11").Value = 333
$ws.Range("E11").Value = 12591
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 8273
$ws.Range("B12").Value = 153548
$ws.Range("C12").Value = 961
$ws.Range("D12").Value = 114990
$ws.Range("E12").Value = 34309
$ws.Range("G12").Value = 27
$ws.Range("H12").Value = 4249
$ws.Range("B17").Value = 80556
$ws.Range("C17").Value = 414
$ws.Range("D17").Value = 41098
$ws.Range("E17").Value = 33395
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 6063
$ws.Range("E40").Value = 5848
$ws.Range("G40").Value = 9
$ws.Range("H40").Value = 1156
$ws.Range("D51").Value = 3032
$ws.Range("E51").Value = 5847
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 404
$ws.Range("B58").Value = 7211
$ws.Range("C58").Value = 78
$ws.Range("D58").Value = 4280
$ws.Range("E58").Value = 2735
$ws.Range("B86").Value = 1908
$ws.Range("C86").Value = 8
$ws.Range("D86").Value = 1603
$ws.Range("E86").Value = 225
$ws.Range("G86").Value = 1
$ws.Range("H86").Value = 80
$ws.Range("B102").Value = 1109
$ws.Range("C102").Value = 80
$ws.Range("D102").Value = 366
$ws.Range("E102").Value = 693
$ws.Range("H102").Value = 50
$ws.Range("B103").Value = 1089
$ws.Range("C103").Value = 0
$ws.Range("D103").Value = 42
$ws.Range("E103").Value = 1041
$ws.Range("H103").Value = 6
$ws.Range("B104").Value = 1064
$ws.Range("C104").Value = 8
$ws.Range("D104").Value = 1029
$ws.Range("E104").Value = 31
$ws.Range("H104").Value = 4
$ws.Range("C105").Value = 17
$ws.Range("D105").Value = 604
$ws.Range("E105").Value = 432
$ws.Range("H105").Value = 9
$ws.Range("B106").Value = 1045
$ws.Range("D106").Value = 862
$ws.Range("E106").Value = 136
$ws.Range("H106").Value = 47
$ws.Range("B137").Value = 417
$ws.Range("C137").Value = 19
$ws.Range("D137").Value = 346
$ws.Range("E137").Value = 69
$ws.Range("B138").Value = 405
$ws.Range("C138").Value = 34
$ws.Range("D138").Value = 131
$ws.Range("E138").Value = 272
$ws.Range("H138").Value = 2
$ws.Range("C139").Value = 9
$ws.Range("D139").Value = 123
$ws.Range("E139").Value = 270
$ws.Range("H139").Value = 5
$ws.Range("D163").Value = 61
$ws.Range("E163").Value = 71
